$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cameras")

# Widen column C slightly
$ws.Columns.Item(3).ColumnWidth = 11.33203125

# New data: numeric ID, Name, Price
$data = @(
    @(1,  "LAMBORGHINI SIAN",     3600000),
    @(2,  "TESLA ROADSTER",        400000),
    @(3,  "RENAULT ALPHINE",       270000),
    @(4,  "CHEVROLET CORVETTE",    110000),
    @(5,  "MERCEDES AMG ONE",     2700000),
    @(6,  "BUGATTI DIVO",         5800000),
    @(7,  "HONDA NSX",             320000),
    @(8,  "SUBARU BRZ",             90000),
    @(9,  "AUDI E-TRON GT",        205000),
    @(10, "KIA STINGER",            95000),
    @(11, "BMW M4",                280000),
    @(12, "BENTLEY BENTAYGA",      710000)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).NumberFormat = "0"
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

$ws.Range("A17").Select()
